$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 689 entirely ("明白なるアラブの言葉によって" post),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(689).Delete()
